# overview.xlsx template update
# - course.full_name -> course.name (used as the placeholder value in the
#   sample/export row, cell D2 on "Kursliste")
# - the saved cursor/selection moves from D10 to D7
# - a new (empty, underline-styled) input cell appears at D7, matching the
#   existing style already used by E7/E9/G9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the placeholder text in D2 from "course.full_name" to "course.name".
$ws.Range("D2").Value = "course.name"

# Give D7 the same "underline" cell style already used by the neighbouring
# input placeholder cells (E7, E9, G9).
$ws.Range("D7").Font.Underline = 2

# Move the active selection to D7 (was D10).
$ws.Range("D7").Select() | Out-Null
